$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A, shifting the existing table right.
$ws.Columns("A").Insert()

# 2. Header cell for the new column.
$ws.Range("A1").Value = "Match ID"
$ws.Range("A1").Font.Bold = $true

# 3. Styled-but-empty cells on the duplicate header row and the blank separator row.
$ws.Range("A2").Font.Bold = $true
$ws.Range("A3").Font.Bold = $true

# 4. Fill the Match ID values for the visible + hidden detail rows (4-19).
$ws.Range("A4:A19").Value = 21
$ws.Range("A4:A19").Font.Bold = $true

# 5. Totals row keeps the default (unbolded) style.
$ws.Range("A20").Value = 21
$ws.Rows(20).AutoFit()

# 6. Update the active selection to match the authored state.
$ws.Range("A1:A19").Select()
